$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 66; $row++) {
    $ws.Cells.Item($row, 6).Value = "Done"
}
